# Update the header date line.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-03-06 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-03-07 Thursday", 2) | Out-Null

# Update each division expression in the practice table.
#
# NOTE: this runtime's Find.Execute always searches/replaces against the
# *whole* document (first match in document order) no matter which Range's
# .Find is invoked - it is not actually scoped to the calling Range. Several
# of the new values here happen to equal *other* cells' old values (e.g.
# "31÷8=" and "24÷3=" each appear as both an old value in one cell and a new
# value in another), so a naive per-cell Find/Replace can clobber the wrong
# cell. To avoid that entirely, address each cell by its table position and
# overwrite its Range.Text directly (all replacement strings are the same
# length as the originals, so cell character offsets stay valid throughout).
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "39÷9="
$cell = $t.Cell(1, 2); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "26÷8="
$cell = $t.Cell(1, 3); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "17÷3="
$cell = $t.Cell(1, 4); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "49÷2="
$cell = $t.Cell(1, 5); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "45÷8="

$cell = $t.Cell(5, 1); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "36÷4="
$cell = $t.Cell(5, 2); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "73÷8="
$cell = $t.Cell(5, 3); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "40÷6="
$cell = $t.Cell(5, 4); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "56÷7="
$cell = $t.Cell(5, 5); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "42÷8="

$cell = $t.Cell(9, 1); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "46÷4="
$cell = $t.Cell(9, 2); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "17÷5="
$cell = $t.Cell(9, 3); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "31÷8="
$cell = $t.Cell(9, 4); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "51÷5="
$cell = $t.Cell(9, 5); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "20÷8="

$cell = $t.Cell(13, 1); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "43÷3="
$cell = $t.Cell(13, 2); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "80÷9="
$cell = $t.Cell(13, 3); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "68÷4="
$cell = $t.Cell(13, 4); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "24÷3="
$cell = $t.Cell(13, 5); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "87÷5="

$cell = $t.Cell(17, 1); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "56÷8="
$cell = $t.Cell(17, 2); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "49÷8="
$cell = $t.Cell(17, 3); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "64÷4="
$cell = $t.Cell(17, 4); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "76÷6="
$cell = $t.Cell(17, 5); $rng = $cell.Range; $d.Range($rng.Start, $rng.End - 1).Text = "72÷9="

Write-Host "Done"
